# Update the roster sheet so that each player's Position (col B) and
# Team (col C) reflect their current assignment. Player names in column A
# keep their existing row order; only columns B and C are corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Russell Westbrook", "PG",    "Denver Nuggets"),
    @("Ayo Dosunmu",        "SG,SF", "Chicago Bulls"),
    @("Chris Paul",         "PG",    "San Antonio Spurs"),
    @("Jalen Green",        "PG,SG", "Houston Rockets"),
    @("Deni Avdija",        "SF,PF", "Portland Trail Blazers"),
    @("Pascal Siakam",      "SF,PF", "Indiana Pacers"),
    @("Naz Reid",           "PF,C",  "Minnesota Timberwolves"),
    @("Jerami Grant",       "SF,PF", "Portland Trail Blazers"),
    @("Rudy Gobert",        "C",     "Minnesota Timberwolves"),
    @("Nikola Jokic",       "C",     "Denver Nuggets"),
    @("Jakob Poeltl",       "C",     "Toronto Raptors"),
    @("Clint Capela",       "C",     "Atlanta Hawks"),
    @("Dejounte Murray",    "PG,SG", "New Orleans Pelicans"),
    @("Jaylen Brown",       "SG,SF", "Boston Celtics"),
    @("Jalen Suggs",        "PG,SG", "Orlando Magic"),
    @("Paolo Banchero",     "SF,PF", "Orlando Magic"),
    @("Chet Holmgren",      "PF,C",  "Oklahoma City Thunder"),
    @("Bogdan Bogdanovic",  "SG,SF", "Atlanta Hawks")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
